# Add a "Save" column (column H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell - set value first, then copy the formatting from an existing
# header cell (G1) so it reuses the same cell style as the rest of the header row.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for the new column (plain, unstyled numeric cells like the rest of the data)
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
